# Append run: new scrape timestamp 2026-02-13 02:02:35 (JST), with two new
# listings inserted into the ranking (row 7: PHP/Laravel; row 11: 放置中の
# 法人ドメイン), pushing the previously-existing rows down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-13 02:02:35"

# --- 1) Insert the two new rows first, so every other row's data just needs
#        its timestamp refreshed in place afterwards. ---------------------

# New row for "PHP/Laravelエンジニア募集..." goes in right after the current
# row 6 ("自動化システム"), i.e. at row 7 - shifts old rows 7-10 to 8-11.
$ws.Rows.Item(7).Insert()

# New row for "放置中の法人ドメインを持っている企業様" goes in right after
# what is now row 10 ("プロジェクトマネジメント", shifted from old row 9),
# i.e. at row 11 - shifts old row 10 (now 11, "SES経営者向け...") to 12.
$ws.Rows.Item(11).Insert()

# --- 2) Remove all existing hyperlinks; row-insert does not re-target them
#        to the shifted cells, so rebuild the hyperlink set from scratch
#        once every row is in its final place. ----------------------------
$ws.Hyperlinks.Delete()

# --- 3) Refresh the "取得日時" timestamp on every data row (2-12). --------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 4) Fill in the new row 7: PHP/Laravelエンジニア募集 ------------------
$ws.Cells.Item(7, 2).Value = "PHP/Laravelエンジニア募集(大規模Webシステム/フルリモート)"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5490679"
$ws.Cells.Item(7, 7).Value = 55
$ws.Cells.Item(7, 8).Value = "○PHP"

# --- 5) Fill in the new row 11: 放置中の法人ドメインを持っている企業様 ----
$ws.Cells.Item(11, 2).Value = "放置中の法人ドメインを持っている企業様"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5490905"
$ws.Cells.Item(11, 7).Value = 10

# --- 6) Rebuild the hyperlinks for column F, rows 2-12, from the URL text
#        now sitting in each cell. -----------------------------------------
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value2)
}

# dimension / used range will naturally extend to H12 once H7 has a value
# and rows up to 12 carry data.
